$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $ws.Range("B4").Style
}

$updates = @(
    @{Cell="D2"; Value="43.204.53"},
    @{Cell="E2"; Value="  -0.03%  "},
    @{Cell="D3"; Value="2.403.74"},
    @{Cell="E3"; Value="  +5.08%  "},
    @{Cell="E4"; Value="  -0.35%  "},
    @{Cell="D5"; Value="334.16"},
    @{Cell="E5"; Value="  +8.44%  "},
    @{Cell="D6"; Value="105.95"},
    @{Cell="E6"; Value="  -6.97%  "},
    @{Cell="E7"; Value="  +2.87%  "},
    @{Cell="E8"; Value="  -0.10%  "},
    @{Cell="D9"; Value="0.652"},
    @{Cell="E9"; Value="  +5.92%  "},
    @{Cell="D10"; Value="42.37"},
    @{Cell="E10"; Value="  -5.55%  "},
    @{Cell="E11"; Value="  +1.30%  "},
    @{Cell="D12"; Value="8.77"},
    @{Cell="E12"; Value="  -1.79%  "},
    @{Cell="E13"; Value="  +0.93%  "},
    @{Cell="D14"; Value="17.26"},
    @{Cell="E14"; Value="  +11.63%  "},
    @{Cell="E15"; Value="  +2.00%  "},
    @{Cell="D16"; Value="2.761.70"},
    @{Cell="E16"; Value="  +5.02%  "},
    @{Cell="D17"; Value="2.395.82"},
    @{Cell="E17"; Value="  +5.02%  "},
    @{Cell="D18"; Value="43.208.94"},
    @{Cell="E18"; Value="  +0.15%  "},
    @{Cell="E19"; Value="  +6.05%  "},
    @{Cell="E20"; Value="  +1.23%  "},
    @{Cell="D21"; Value="3.85"},
    @{Cell="E21"; Value="  +6.52%  "},
    @{Cell="D22"; Value="77.24"},
    @{Cell="E22"; Value="  +2.35%  "},
    @{Cell="D23"; Value="277.22"},
    @{Cell="E23"; Value="  +8.97%  "},
    @{Cell="D24"; Value="2.43"},
    @{Cell="E24"; Value="  -1.84%  "},
    @{Cell="D25"; Value="9.87"},
    @{Cell="E25"; Value="  +9.24%  "},
    @{Cell="D26"; Value="11.98"},
    @{Cell="E26"; Value="  +1.72%  "},
    @{Cell="E27"; Value="  -0.03%  "},
    @{Cell="D28"; Value="23.29"},
    @{Cell="E28"; Value="  +4.36%  "},
    @{Cell="E29"; Value="  -1.81%  "},
    @{Cell="D30"; Value="175.98"},
    @{Cell="E30"; Value="  +0.62%  "},
    @{Cell="D31"; Value="37.30"},
    @{Cell="E31"; Value="  -2.84%  "},
    @{Cell="B32"; Value="Hedera"},
    @{Cell="C32"; Value="https://coinranking.com/coin/jad286TjB+hedera-hbar"},
    @{Cell="D32"; Value="0.0942"},
    @{Cell="E32"; Value="  +4.18%  "},
    @{Cell="B33"; Value="WEMIXToken"},
    @{Cell="C33"; Value="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"},
    @{Cell="D33"; Value="3.15"},
    @{Cell="E33"; Value="  -0.93%  "},
    @{Cell="D34"; Value="6.12"},
    @{Cell="E34"; Value="  +7.35%  "},
    @{Cell="E35"; Value="  +5.28%  "},
    @{Cell="E36"; Value="  -3.68%  "},
    @{Cell="D37"; Value="4.12"},
    @{Cell="E38"; Value="  -3.43%  "},
    @{Cell="D39"; Value="0.109"},
    @{Cell="E39"; Value="  +3.03%  "},
    @{Cell="E40"; Value="  +11.05%  "},
    @{Cell="E41"; Value="  +12.18%  "},
    @{Cell="D42"; Value="0.237"},
    @{Cell="E42"; Value="  +1.39%  "},
    @{Cell="D43"; Value="70.50"},
    @{Cell="E43"; Value="  -2.95%  "},
    @{Cell="D44"; Value="122.16"},
    @{Cell="E44"; Value="  +13.05%  "},
    @{Cell="E45"; Value="  +0.10%  "},
    @{Cell="D46"; Value="92.29"},
    @{Cell="E46"; Value="  +43.40%  "},
    @{Cell="E47"; Value="  -2.67%  "},
    @{Cell="E48"; Value="  -1.23%  "},
    @{Cell="D49"; Value="9.29"},
    @{Cell="E49"; Value="  +5.24%  "},
    @{Cell="D50"; Value="0.509"},
    @{Cell="E50"; Value="  +15.76%  "},
    @{Cell="E51"; Value="  +1.22%  "}
)

foreach ($u in $updates) {
    Set-TextValue $ws $u.Cell $u.Value
}
